# Applies the "complete: n = N" -> "n = N" cleanup across all three tables,
# plus a handful of variable-label wording tweaks (Table 2 / Table 3).

$wb = $excel.ActiveWorkbook

# 1) Strip the leading "complete: " from every "complete: n = N" line,
#    wherever it occurs, in every worksheet of the workbook.
foreach ($ws in $wb.Worksheets) {
    $used = $ws.UsedRange
    $rowCount = $used.Rows.Count
    $colCount = $used.Columns.Count
    $rowStart = $used.Row
    $colStart = $used.Column

    for ($r = 0; $r -lt $rowCount; $r++) {
        for ($c = 0; $c -lt $colCount; $c++) {
            $cell = $ws.Cells.Item($rowStart + $r, $colStart + $c)
            $v = $cell.Value2
            if ($v -ne $null -and $v -is [string] -and $v -like "*complete: n = *") {
                $cell.Value2 = $v -replace "complete: n = ", "n = "
            }
        }
    }
}

# 2) Variable-label wording changes.

# Table 2
$ws2 = $wb.Worksheets.Item("Table 2")
$ws2.Range("A4").Value2 = "CT abnormality (CT score ≥ 1)"

# Table 3
$ws3 = $wb.Worksheets.Item("Table 3")
$ws3.Range("A4").Value2  = "SMWD < ref."
$ws3.Range("A5").Value2  = "Fatigue score (likert CFS)"
$ws3.Range("A7").Value2  = "General health score (EQ5D5L VAS)"
$ws3.Range("A8").Value2  = "Imp. general health (VAS < 73, EQ5D5L)"
$ws3.Range("A9").Value2  = "Mobility impairment score (EQ5D5L)"
$ws3.Range("A10").Value2 = "Imp. mobility (score  > 1, EQ5D5L)"
$ws3.Range("A11").Value2 = "Self-care impairment score (EQ5D5L)"
$ws3.Range("A12").Value2 = "Imp. self-care (score  > 1, EQ5D5L)"
$ws3.Range("A13").Value2 = "Activity impairment score (EQ5D5L)"
$ws3.Range("A14").Value2 = "Imp. usual activity (score  > 1, EQ5D5L)"
$ws3.Range("A15").Value2 = "Pain/discomfort score (EQ5D5L)"
$ws3.Range("A16").Value2 = "Pain/discomfort (score  > 1, EQ5D5L)"
$ws3.Range("A17").Value2 = "Anxiety/depression score (EQ5D5L)"
$ws3.Range("A18").Value2 = "Anxiety/depression (score  > 1, EQ5D5L)"
$ws3.Range("A19").Value2 = "Stress score (PSS)"
$ws3.Range("A21").Value2 = "Somatic symptom disorder score (SSD-12)"
$ws3.Range("A22").Value2 = "Resilience score (BRCS)"
